# Adds a new handback record (d08ad933-095b-40df-b7c4-daa762fa9ddc) as row 4
# to the Overview, zh-cn and de-de sheets, mirroring the existing
# "in sync with en-US" rows (the 0dea2988... entries).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Cells.Item(4, 1).Value = "d08ad933-095b-40df-b7c4-daa762fa9ddc.md"
$ov.Cells.Item(4, 2).Value = "e2e\d08ad933-095b-40df-b7c4-daa762fa9ddc.md"
$ov.Cells.Item(4, 3).Value = ".md"
$ov.Cells.Item(4, 5).Value = "Handed back: in sync with en-US"
$ov.Cells.Item(4, 6).Value = "Handed back: in sync with en-US"
$ov.Cells.Item(4, 7).Value = "2016-08-18 12:45:15"
$ov.Cells.Item(4, 7).NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ov.Hyperlinks.Add($ov.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/master/e2e/d08ad933-095b-40df-b7c4-daa762fa9ddc.md", "", "", "e2e\d08ad933-095b-40df-b7c4-daa762fa9ddc.md")
$ov.Range("B4").Style = "HyperLink"

$ovTable = $ov.ListObjects.Item(1)
$ovTable.Resize($ov.Range("A1:G4"))

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Cells.Item(4, 1).Value = "d08ad933-095b-40df-b7c4-daa762fa9ddc.md"
$zh.Cells.Item(4, 2).Value = ".md"
$zh.Cells.Item(4, 3).Value = "Handed back: in sync with en-US"
$zh.Cells.Item(4, 4).Value = "e2e"
$zh.Cells.Item(4, 5).Value = "ht"
$zh.Cells.Item(4, 6).Value = "True"
$zh.Cells.Item(4, 7).Value = "d08ad933-095b-40df-b7c4-daa762fa9ddc.bba7348ffe03113c13ca99620c42166a1a526839.zh-cn.xlf"
$zh.Cells.Item(4, 8).Value = "2016-08-18 12:45:06"
$zh.Cells.Item(4, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zh.Cells.Item(4, 9).Value = "d08ad933-095b-40df-b7c4-daa762fa9ddc.md"
$zh.Cells.Item(4, 10).Value = "d08ad933-095b-40df-b7c4-daa762fa9ddc.bba7348ffe03113c13ca99620c42166a1a526839.zh-cn.xlf"
$zh.Cells.Item(4, 11).Value = "2016-08-18 12:45:30"
$zh.Cells.Item(4, 11).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zh.Cells.Item(4, 12).Value = ""
$zh.Cells.Item(4, 13).Value = "True"
$zh.Cells.Item(4, 14).Value = ""
$zh.Cells.Item(4, 15).Value = "False"
$zh.Cells.Item(4, 16).Value = ""

$zh.Hyperlinks.Add($zh.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/master/e2e/d08ad933-095b-40df-b7c4-daa762fa9ddc.md", "", "", "d08ad933-095b-40df-b7c4-daa762fa9ddc.md")
$zh.Range("A4").Style = "HyperLink"
$zh.Hyperlinks.Add($zh.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/master/e2e/d08ad933-095b-40df-b7c4-daa762fa9ddc.md", "", "", "d08ad933-095b-40df-b7c4-daa762fa9ddc.md")
$zh.Range("I4").Style = "HyperLink"

$zhTable = $zh.ListObjects.Item(1)
$zhTable.Resize($zh.Range("A1:P4"))

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Cells.Item(4, 1).Value = "d08ad933-095b-40df-b7c4-daa762fa9ddc.md"
$de.Cells.Item(4, 2).Value = ".md"
$de.Cells.Item(4, 3).Value = "Handed back: in sync with en-US"
$de.Cells.Item(4, 4).Value = "e2e"
$de.Cells.Item(4, 5).Value = "ht"
$de.Cells.Item(4, 6).Value = "True"
$de.Cells.Item(4, 7).Value = "d08ad933-095b-40df-b7c4-daa762fa9ddc.bba7348ffe03113c13ca99620c42166a1a526839.de-de.xlf"
$de.Cells.Item(4, 8).Value = "2016-08-18 12:45:15"
$de.Cells.Item(4, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$de.Cells.Item(4, 9).Value = "d08ad933-095b-40df-b7c4-daa762fa9ddc.md"
$de.Cells.Item(4, 10).Value = "d08ad933-095b-40df-b7c4-daa762fa9ddc.bba7348ffe03113c13ca99620c42166a1a526839.de-de.xlf"
$de.Cells.Item(4, 11).Value = "2016-08-18 12:45:38"
$de.Cells.Item(4, 11).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$de.Cells.Item(4, 12).Value = ""
$de.Cells.Item(4, 13).Value = "True"
$de.Cells.Item(4, 14).Value = ""
$de.Cells.Item(4, 15).Value = "False"
$de.Cells.Item(4, 16).Value = ""

$de.Hyperlinks.Add($de.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/master/e2e/d08ad933-095b-40df-b7c4-daa762fa9ddc.md", "", "", "d08ad933-095b-40df-b7c4-daa762fa9ddc.md")
$de.Range("A4").Style = "HyperLink"
$de.Hyperlinks.Add($de.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/master/e2e/d08ad933-095b-40df-b7c4-daa762fa9ddc.md", "", "", "d08ad933-095b-40df-b7c4-daa762fa9ddc.md")
$de.Range("I4").Style = "HyperLink"

$deTable = $de.ListObjects.Item(1)
$deTable.Resize($de.Range("A1:P4"))

Write-Host "Handback row added."
